# Rename the "AT commands" shared strings used on the Command sheet from
# human-readable labels to UPPER_SNAKE_CASE command identifiers, and change
# which sheet/cell is active so the workbook reopens on the "Command" sheet
# with E25 selected (scrolled so row 10 is at top).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Command")

# --- Rename the AT command labels in column E/F/G/H ---
$ws.Range("E2").Value   = "ENTER"
$ws.Range("E3:E4").Value = "TAB_AND_SHIFT_TAB"
$ws.Range("F5").Value   = "TAB_AND_SHIFT_TAB"
$ws.Range("E5").Value   = "F_AND_SHIFT_F"
$ws.Range("G5").Value   = "UP_AND_DOWN"
$ws.Range("E6:E9").Value = "UP_AND_DOWN"
$ws.Range("H5").Value   = "LEFT_AND_RIGHT (with Smart Navigation on)"
$ws.Range("F7").Value   = "S (Navigate by first letter of menuitem)"
$ws.Range("F9").Value   = "B (Navigate by first letter of menuitem)"
$ws.Range("E10:E13").Value = "ENTER"
$ws.Range("F10:F13").Value = "SPACE"
$ws.Range("E14:E15").Value = "ESC"
$ws.Range("E16:E25").Value = "INSERT_TAB"
$ws.Range("F16:F25").Value = "INSERT_UP"

# --- Switch the active sheet/selection ---
# Previously "Tasks" (index 0) was the tab shown on open; now it's "Command".
$ws.Activate()

# Scroll so row 10 is near the top of the view, then land the selection on E25
# (best-effort: the interop layer may not persist the scroll position itself,
# but the active sheet + selected cell below are what matters for the saved view).
$win = $excel.ActiveWindow
$win.ScrollRow = 10
$win.ScrollColumn = 2

$ws.Range("E25").Select()
